$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing unique_id values in column J for rows 3, 5, 7, 9
$ws.Range("J3").Value = "VEC-015-03-193"
$ws.Range("J5").Value = "VEC-015-04-192"
$ws.Range("J7").Value = "VEC-015-04-201"
$ws.Range("J9").Value = "VEC-015-04-200"

# Update selection to J9 to match the final saved selection state
$ws.Range("J9").Select()
